# Update MenuOptions sheet (ebay tickets & experiences menu) and reset the
# active selection to A8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MenuOptions")

$ws.Range("A1").Value = "Tickets & Experiences"
$ws.Range("A2").Value = "Concert Tickets"
$ws.Range("A3").Value = "Sports Tickets"
$ws.Range("A4").Value = "Theater Tickets"
$ws.Range("A5").Value = "Theme Park & Club Passes"
$ws.Range("A6").Value = "Parking Passes"
$ws.Range("A7").Value = "Special Experiences"
$ws.Range("A8").Value = "Other Tickets & Experiences"

$ws.Activate()
$ws.Range("A8").Select()
